# Add nicer phone numbers to data
# - introduces a new phone number "858-222-1234"
# - rotates/varies the PhoneNumber values for several leads (both the
#   "leads as columns" sheet (Sheet1) and the "leads as rows" sheet (Sheet2))
# - bumps the LeadScore for the first three leads (Jennifer, Mark, Saurabh)
#   from 45 to 97 on both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1: fields as rows, leads as columns B..P
$ws2 = $wb.Worksheets.Item(2)   # Sheet2: leads as rows, fields as columns

$newPhone = "858-222-1234"

# --- Sheet1: row 4 is "PhoneNumber", columns B..P are the 15 leads ---
$ws1.Range("E4").Value = $newPhone
$ws1.Range("G4").Value = "212-555-1234"
$ws1.Range("H4").Value = "212-555-4321"
$ws1.Range("I4").Value = "212-555-9999"
$ws1.Range("J4").Value = $newPhone
$ws1.Range("K4").Value = "301-555-4321"
$ws1.Range("L4").Value = "212-555-1234"
$ws1.Range("M4").Value = "212-555-4321"
$ws1.Range("N4").Value = "212-555-9999"
$ws1.Range("O4").Value = $newPhone
$ws1.Range("P4").Value = "301-555-4321"

# --- Sheet1: row 17 is "LeadScore" -- first three leads go from 45 to 97 ---
$ws1.Range("B17").Value = 97
$ws1.Range("C17").Value = 97
$ws1.Range("D17").Value = 97

# --- Sheet2: column D is "PhoneNumber", one lead per row (rows 2..16) ---
$ws2.Range("D5").Value = $newPhone
$ws2.Range("D7").Value = "212-555-1234"
$ws2.Range("D8").Value = "212-555-4321"
$ws2.Range("D9").Value = "212-555-9999"
$ws2.Range("D10").Value = $newPhone
$ws2.Range("D11").Value = "301-555-4321"
$ws2.Range("D12").Value = "212-555-1234"
$ws2.Range("D13").Value = "212-555-4321"
$ws2.Range("D14").Value = "212-555-9999"
$ws2.Range("D15").Value = $newPhone
$ws2.Range("D16").Value = "301-555-4321"

# --- Sheet2: column Q is "LeadScore" -- first three leads go from 45 to 97 ---
$ws2.Range("Q2").Value = 97
$ws2.Range("Q3").Value = 97
$ws2.Range("Q4").Value = 97
